# Refined Evaluation to be more exact
# - Insert two new columns (O, P) for "Correct Pred Predicates Parents" /
#   "Correct Pred Predicates Related", shifting the old Object/Entity/Result
#   columns two positions to the right (O->Q, P->R, Q->S, R->T, S->U, T->V, U->W).
# - Rename the existing M/N headers.
# - Populate the new O/P columns with the same counts that used to live in
#   M/N (the "Correct Extracted Predicates with Parents/Related" numbers),
#   since those are now tracked as separate "detected" vs "correct" metrics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank columns at O:P; this shifts existing O:U -> Q:W,
# carrying over values, formatting and styles automatically.
$ws.Range("O1:P1").EntireColumn.Insert()

# Rename the predicate-related headers.
$ws.Range("M1").Value = "Detected Predicates Doc Parent"
$ws.Range("N1").Value = "Detected Predicates Doc Related"

# New header labels for the freshly-inserted columns.
$ws.Range("O1").Value = "Correct Pred Predicates Parents"
$ws.Range("P1").Value = "Correct Pred Predicates Related"

# Fill the new columns' data rows with the same values that M/N already hold.
$ws.Range("O2").Value = $ws.Range("M2").Value2
$ws.Range("P2").Value = $ws.Range("N2").Value2

$ws.Range("O3").Value = $ws.Range("M3").Value2
$ws.Range("P3").Value = $ws.Range("N3").Value2

$ws.Range("O4").Value = $ws.Range("M4").Value2
$ws.Range("P4").Value = $ws.Range("N4").Value2

$ws.Range("O5").Value = $ws.Range("M5").Value2
$ws.Range("P5").Value = $ws.Range("N5").Value2

$ws.Range("O6").Value = $ws.Range("M6").Value2
$ws.Range("P6").Value = $ws.Range("N6").Value2
